# Auto-generated Excel COM-interop edit script
# Applies numeric cell updates (and a handful of cell additions/removals)
# to the Ixion Profits sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 190477600
$ws.Range("I106").Value = 55557200
$ws.Range("J106").Value = 1000000000
$ws.Range("K106").Value = 55557200
$ws.Range("L106").Value = 1000000000
$ws.Range("M106").Value = -55556569
$ws.Range("N106").Value = -1000001262
$ws.Range("H107").Value = 25000700
$ws.Range("I107").Value = 41667150
$ws.Range("J107").Value = 1025
$ws.Range("K107").Value = 41667150
$ws.Range("L107").Value = 1025
$ws.Range("M107").Value = -41665230
$ws.Range("N107").Value = -4865
$ws.Range("H111").Value = 1000001
$ws.Range("I111").Value = 0
$ws.Range("J111").Value = 1000001
$ws.Range("K111").Value = 0
$ws.Range("L111").Value = 3000003
$ws.Range("M111").Value = ""
$ws.Range("N111").Value = -3006137
$ws.Range("H127").Value = 1825.1471
$ws.Range("I127").Value = 701.5
$ws.Range("J127").Value = 2065.9285
$ws.Range("K127").Value = 2104.5
$ws.Range("L127").Value = 6197.7855
$ws.Range("M127").Value = 2855.5
$ws.Range("N127").Value = -16117.7855
$ws.Range("H137").Value = 1926.8889
$ws.Range("I137").Value = 1592.75
$ws.Range("J137").Value = 4600
$ws.Range("K137").Value = 4778.25
$ws.Range("L137").Value = 13800
$ws.Range("M137").Value = -2228.25
$ws.Range("N137").Value = -18900

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4592.222
$ws.Range("I32").Value = 3395.1
$ws.Range("K32").Value = 3395.1
$ws.Range("M32").Value = -3108.1
$ws.Range("H45").Value = 5466.68
$ws.Range("I45").Value = 8785.77
$ws.Range("K45").Value = 8785.77
$ws.Range("M45").Value = -8408.77
$ws.Range("H61").Value = 5519.769
$ws.Range("I61").Value = 6616.524
$ws.Range("K61").Value = 6616.524
$ws.Range("M61").Value = -6404.524
$ws.Range("H74").Value = 1515.973
$ws.Range("I74").Value = 1504.1384
$ws.Range("J74").Value = 1601.4445
$ws.Range("K74").Value = 1504.1384
$ws.Range("L74").Value = 1601.4445
$ws.Range("M74").Value = -630.1384
$ws.Range("N74").Value = -3349.4445
$ws.Range("H77").Value = 1515.973
$ws.Range("I77").Value = 1504.1384
$ws.Range("J77").Value = 1601.4445
$ws.Range("K77").Value = 7520.692
$ws.Range("L77").Value = 8007.2225
$ws.Range("M77").Value = -3152.692
$ws.Range("N77").Value = -16743.2225
$ws.Range("H102").Value = 4631761
$ws.Range("I102").Value = 5293012.5
$ws.Range("J102").Value = 3000
$ws.Range("K102").Value = 5293012.5
$ws.Range("L102").Value = 3000
$ws.Range("M102").Value = -5291390.5
$ws.Range("N102").Value = -6244
$ws.Range("H132").Value = 3311.718
$ws.Range("I132").Value = 1736
$ws.Range("J132").Value = 5350.8823
$ws.Range("K132").Value = 5208
$ws.Range("L132").Value = 16052.6469
$ws.Range("M132").Value = -2678
$ws.Range("N132").Value = -21112.6469
$ws.Range("H136").Value = 5519.769
$ws.Range("I136").Value = 6616.524
$ws.Range("K136").Value = 19849.572
$ws.Range("M136").Value = -17299.572

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1508.16
$ws.Range("I94").Value = 958.4706
$ws.Range("J94").Value = 2676.25
$ws.Range("K94").Value = 958.4706
$ws.Range("L94").Value = 2676.25
$ws.Range("M94").Value = -507.4706
$ws.Range("N94").Value = -3578.25
$ws.Range("H99").Value = 83334590
$ws.Range("I99").Value = 90910190
$ws.Range("K99").Value = 90910190
$ws.Range("M99").Value = -90908692
$ws.Range("H107").Value = 899.75
$ws.Range("I107").Value = 841.8125
$ws.Range("J107").Value = 1131.5
$ws.Range("K107").Value = 841.8125
$ws.Range("L107").Value = 1131.5
$ws.Range("M107").Value = 1078.1875
$ws.Range("N107").Value = -4971.5
$ws.Range("H134").Value = 4106.511
$ws.Range("I134").Value = 4734.8438
$ws.Range("K134").Value = 14204.5314
$ws.Range("M134").Value = -11669.5314

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4680.4614
$ws.Range("I31").Value = 947.4074000000001
$ws.Range("J31").Value = 13079.833
$ws.Range("K31").Value = 947.4074000000001
$ws.Range("L31").Value = 13079.833
$ws.Range("M31").Value = -652.4074000000001
$ws.Range("N31").Value = -13669.833
$ws.Range("H34").Value = 4680.4614
$ws.Range("I34").Value = 947.4074000000001
$ws.Range("J34").Value = 13079.833
$ws.Range("K34").Value = 947.4074000000001
$ws.Range("L34").Value = 13079.833
$ws.Range("M34").Value = -745.4074000000001
$ws.Range("N34").Value = -13483.833
$ws.Range("H107").Value = 282.4091
$ws.Range("I107").Value = 128.83333
$ws.Range("J107").Value = 340
$ws.Range("K107").Value = 128.83333
$ws.Range("L107").Value = 340
$ws.Range("M107").Value = 1791.16667
$ws.Range("N107").Value = -4180
$ws.Range("H132").Value = 1994.8148
$ws.Range("I132").Value = 1782
$ws.Range("J132").Value = 3218.5
$ws.Range("K132").Value = 5346
$ws.Range("L132").Value = 9655.5
$ws.Range("M132").Value = -2816
$ws.Range("N132").Value = -14715.5
$ws.Range("H134").Value = 1547.0741
$ws.Range("I134").Value = 1503.2632
$ws.Range("J134").Value = 1651.125
$ws.Range("K134").Value = 4509.7896
$ws.Range("L134").Value = 4953.375
$ws.Range("M134").Value = -1974.7896
$ws.Range("N134").Value = -10023.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 8484.406999999999
$ws.Range("I3").Value = 11589.917
$ws.Range("J3").Value = 6000
$ws.Range("K3").Value = 34769.751
$ws.Range("L3").Value = 18000
$ws.Range("M3").Value = -34657.751
$ws.Range("N3").Value = -18224

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").Value = ""
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").Value = ""
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").Value = ""
$ws.Range("H102").Value = 1039.1666
$ws.Range("I102").Value = 887
$ws.Range("K102").Value = 887
$ws.Range("M102").Value = 735
$ws.Range("H107").Value = 1974.75
$ws.Range("I107").Value = 977.8
$ws.Range("J107").Value = 2307.0667
$ws.Range("K107").Value = 977.8
$ws.Range("L107").Value = 2307.0667
$ws.Range("M107").Value = 942.2
$ws.Range("N107").Value = -6147.066699999999
$ws.Range("H132").Value = 3938.1667
$ws.Range("I132").Value = 4276.6665
$ws.Range("J132").Value = 3599.6667
$ws.Range("K132").Value = 12829.9995
$ws.Range("L132").Value = 10799.0001
$ws.Range("M132").Value = -10299.9995
$ws.Range("N132").Value = -15859.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1649.875
$ws.Range("I16").Value = 800
$ws.Range("J16").Value = 2499.75
$ws.Range("K16").Value = 800
$ws.Range("L16").Value = 2499.75
$ws.Range("M16").Value = -630
$ws.Range("N16").Value = -2839.75
$ws.Range("H40").Value = 125003070
$ws.Range("I40").Value = 125003070
$ws.Range("K40").Value = 125003070
$ws.Range("M40").Value = -125002934
$ws.Range("H93").Value = 100040560
$ws.Range("I93").Value = 100400
$ws.Range("K93").Value = 100400
$ws.Range("M93").Value = -99152
$ws.Range("H122").Value = 71428570
$ws.Range("I122").Value = 71428570
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 214285710
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -214283260
$ws.Range("N122").Value = ""
$ws.Range("H132").Value = 13984873
$ws.Range("I132").Value = 16673354
$ws.Range("J132").Value = 4772.6
$ws.Range("K132").Value = 50020062
$ws.Range("L132").Value = 14317.8
$ws.Range("M132").Value = -50017532
$ws.Range("N132").Value = -19377.8
$ws.Range("H136").Value = 4453.115
$ws.Range("I136").Value = 3881.5386
$ws.Range("K136").Value = 11644.6158
$ws.Range("M136").Value = -9094.6158

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").Value = ""
$ws.Range("H113").Value = 980.2353000000001
$ws.Range("I113").Value = 1147.3334
$ws.Range("J113").Value = 792.25
$ws.Range("K113").Value = 3442.0002
$ws.Range("L113").Value = 2376.75
$ws.Range("M113").Value = -1272.0002
$ws.Range("N113").Value = -6716.75
$ws.Range("H122").Value = 2414.2856
$ws.Range("I122").Value = 2333.3333
$ws.Range("J122").Value = 2475
$ws.Range("K122").Value = 6999.999899999999
$ws.Range("L122").Value = 7425
$ws.Range("M122").Value = -4549.999899999999
$ws.Range("N122").Value = -12325
$ws.Range("H126").Value = 1349.1666
$ws.Range("I126").Value = 998.3333
$ws.Range("J126").Value = 1700
$ws.Range("K126").Value = 2994.9999
$ws.Range("L126").Value = 5100
$ws.Range("M126").Value = -524.9998999999998
$ws.Range("N126").Value = -10040
$ws.Range("H132").Value = 1890.2727
$ws.Range("I132").Value = 1061.8334
$ws.Range("J132").Value = 2884.4
$ws.Range("K132").Value = 3185.5002
$ws.Range("L132").Value = 8653.200000000001
$ws.Range("M132").Value = -655.5001999999999
$ws.Range("N132").Value = -13713.2
$ws.Range("H136").Value = 2718.3547
$ws.Range("I136").Value = 2889.0454
$ws.Range("J136").Value = 2301.111
$ws.Range("K136").Value = 8667.136200000001
$ws.Range("L136").Value = 6903.333
$ws.Range("M136").Value = -6117.136200000001
$ws.Range("N136").Value = -12003.333
